# Update "lines_states" sheet with refreshed contingency data (rene fine)
# and append two new contingency rows (line7 / line8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-15: update from_bus/to_bus/in_service values; rows 8-9 are renamed
# from "extrN" to "lineN" to make room for the new line7/line8 entries, with
# the remaining extr rows shifting names down by two (extr3->extr1, ...).
$rowData = @(
    @{ Row = 8;  Name = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  Name = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; Name = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; Name = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; Name = "extr3"; C = 10; D = 11; E = $true  },
    @{ Row = 13; Name = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; Name = "extr5"; C = 9;  D = 11; E = $true  },
    @{ Row = 15; Name = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; Name = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; Name = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($item in $rowData) {
    $r = $item.Row

    $ws.Cells.Item($r, 2).Value = $item.Name
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}

# The two brand-new rows (16 and 17) need column A filled in (sequential
# index continuing from row 15) and the bordered/bold style copied over from
# the row above so they match the rest of the table.
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(17, 1).Value = 15

$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
